$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- MOVIE_DETAILS table (columns G:H, rows 13-18): remove the
# "DESCRIPTION / VARCHAR(255)" row (row 14) and shift the rows below it
# up by one, leaving row 18 blank (then cleared below).
for ($r = 14; $r -le 17; $r++) {
    $nr = $r + 1
    $ws.Range("G$r").Value = $ws.Range("G$nr").Value2
    $ws.Range("H$r").Value = $ws.Range("H$nr").Value2
}
for ($r = 14; $r -le 17; $r++) {
    $nr = $r + 1
    $ws.Range("G$nr`:H$nr").Copy()
    $ws.Range("G$r`:H$r").PasteSpecial(-4122)
}
$ws.Application.CutCopyMode = $false
$ws.Range("G18:H18").Clear()

# --- GENRE / VARCHAR(255) row (D17:E17) gets highlighted the same yellow
# used by the table headers (fill colour only, no centered alignment), so
# copy the border from a plain bordered cell first and then change the
# fill colour to the header yellow.
$ws.Range("D5").Copy()
$ws.Range("D17:E17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D17:E17").Interior.Color = 65535

# --- Selection moves to O14 (reflecting wherever the user clicked next).
$ws.Range("O14").Select()
